# Apply the "recorder-facility" StructureDefinition 5.0.0 -> 6.0.0 update
# (Alvearie FHIR IG publish for commit 8e4a450c507ef6f746e072652acbb72e9504f19a)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Metadata": version bump, new date, publisher + jurisdiction
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Row 11 duplicated the "Contact" / "No display for ContactDetail" row
# that was already present as row 10; drop the duplicate so the
# "Jurisdiction" row introduced below lines the table back up with the
# rows that follow (Description, Purpose, Copyright, ...).
$meta.Rows.Item(11).Delete()

# Version 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date refreshed to the new publication timestamp
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher now has a real value
$meta.Range("B9").Value = "Alvearie Team"

# The old duplicate "Contact" row (row 10) becomes the new
# "Jurisdiction" / "United States of America" row
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# ---------------------------------------------------------------------
# Sheet "Elements": give the root Extension row its real short/definition
# ---------------------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

$elements.Range("K2").Value = "Recorder Facility"
$elements.Range("L2").Value = "Facility where condition was recorded"
